$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.987.75"
$ws.Range("E2").Value = "  -0.24%  "

$ws.Range("D3").Value = "1.676.97"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.519"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.45%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  +0.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.66%  "

$ws.Range("E11").Value = "  -0.30%  "

$ws.Range("D12").Value = "1.913.12"
$ws.Range("E12").Value = "  +0.23%  "

$ws.Range("D13").Value = "1.684.64"
$ws.Range("E13").Value = "  +0.56%  "

$ws.Range("E14").Value = "  -0.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.529"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.45%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.06%  "

$ws.Range("D17").Value = "26.985.03"
$ws.Range("E17").Value = "  -0.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "237.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.23%  "

$ws.Range("E20").Value = "  -0.54%  "

$ws.Range("E21").Value = "  +0.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.66%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.82%  "

$ws.Range("E24").Value = "  -1.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.85%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.41%  "

$ws.Range("E28").Value = "  -1.50%  "

$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("E31").Value = "  -0.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.06%  "

$ws.Range("D33").Value = "1.479.85"
$ws.Range("E33").Value = "  +1.73%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.69"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.87%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.584"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0173"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.903"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.84%  "

$ws.Range("E40").Value = "  -3.31%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.03"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.38%  "

$ws.Range("E42").Value = "  +0.06%  "

$ws.Range("E43").Value = "  +2.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "67.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.01%  "

$ws.Range("D45").Value = "1.820.13"
$ws.Range("E45").Value = "  +0.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.781"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.28%  "

$ws.Range("E48").Value = "  +1.53%  "

$ws.Range("E49").Value = "  -0.61%  "

$ws.Range("E50").Value = "  +2.09%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.39%  "
